$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J39").Value = 1483.3334
$ws.Range("L39").Value = 4450.0002
$ws.Range("N39").Value = -5042.0002
$ws.Range("H39").Value = 556.7778
$ws.Range("K39").Value = 280.5
$ws.Range("M39").Value = 15.5
$ws.Range("I39").Value = 93.5
$ws.Range("K42").Value = 246.85713
$ws.Range("M42").Value = -16.85712999999998
$ws.Range("H42").Value = 97
$ws.Range("I42").Value = 82.28570999999999
$ws.Range("H48").Value = 0
$ws.Range("N48").Value = $null
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null
$ws.Range("M61").Value = -210.875
$ws.Range("I61").Value = 127.625
$ws.Range("H61").Value = 127.625
$ws.Range("K61").Value = 382.875
$ws.Range("I112").Value = 166667180
$ws.Range("L112").Value = 7478.3334
$ws.Range("M112").Value = -500000432
$ws.Range("J112").Value = 2492.7778
$ws.Range("K112").Value = 500001540
$ws.Range("H112").Value = 41668664
$ws.Range("N112").Value = -9694.3334
$ws.Range("L137").Value = 16312.6362
$ws.Range("K137").Value = 6982946.399999999
$ws.Range("H137").Value = 1854605.9
$ws.Range("N137").Value = -21412.6362
$ws.Range("M137").Value = -6980396.399999999
$ws.Range("J137").Value = 5437.5454
$ws.Range("I137").Value = 2327648.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J32").Value = 23184
$ws.Range("L32").Value = 23184
$ws.Range("N32").Value = -23758
$ws.Range("I32").Value = 11186.418
$ws.Range("K32").Value = 11186.418
$ws.Range("M32").Value = -10899.418
$ws.Range("H32").Value = 15145.62
$ws.Range("J61").Value = 112866.664
$ws.Range("M61").Value = -41709572
$ws.Range("I61").Value = 41709784
$ws.Range("H61").Value = 30365172
$ws.Range("N61").Value = -113290.664
$ws.Range("K61").Value = 41709784
$ws.Range("L61").Value = 112866.664
$ws.Range("L74").Value = 96790.91
$ws.Range("I74").Value = 6967168
$ws.Range("J74").Value = 96790.91
$ws.Range("K74").Value = 6967168
$ws.Range("N74").Value = -98538.91
$ws.Range("H74").Value = 5686250.5
$ws.Range("M74").Value = -6966294
$ws.Range("H77").Value = 5686250.5
$ws.Range("J77").Value = 96790.91
$ws.Range("N77").Value = -492690.55
$ws.Range("M77").Value = -34831472
$ws.Range("L77").Value = 483954.55
$ws.Range("I77").Value = 6967168
$ws.Range("K77").Value = 34835840
$ws.Range("N110").Value = -6336.4
$ws.Range("M110").Value = -712669.6
$ws.Range("L110").Value = 2246.4
$ws.Range("J110").Value = 2246.4
$ws.Range("I110").Value = 714714.6
$ws.Range("K110").Value = 714714.6
$ws.Range("H110").Value = 527223
$ws.Range("H119").Value = 30500
$ws.Range("J119").Value = 30500
$ws.Range("L119").Value = 30500
$ws.Range("N119").Value = -40176
$ws.Range("H122").Value = 4446664.5
$ws.Range("I122").Value = 2083.389
$ws.Range("J122").Value = 15875587
$ws.Range("M122").Value = -3800.167
$ws.Range("N122").Value = -47631661
$ws.Range("K122").Value = 6250.167
$ws.Range("L122").Value = 47626761
$ws.Range("J132").Value = 139484.75
$ws.Range("L132").Value = 418454.25
$ws.Range("H132").Value = 10684755
$ws.Range("N132").Value = -423514.25
$ws.Range("I132").Value = 12847887
$ws.Range("M132").Value = -38541131
$ws.Range("K132").Value = 38543661
$ws.Range("N136").Value = -343699.992
$ws.Range("M136").Value = -125126802
$ws.Range("L136").Value = 338599.992
$ws.Range("H136").Value = 30365172
$ws.Range("J136").Value = 112866.664
$ws.Range("K136").Value = 125129352
$ws.Range("I136").Value = 41709784

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K7").Value = 2950
$ws.Range("I7").Value = 2950
$ws.Range("M7").Value = -2837
$ws.Range("N7").Value = -30026
$ws.Range("H7").Value = 8320
$ws.Range("L7").Value = 29800
$ws.Range("J7").Value = 29800

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 52507.95
$ws.Range("H31").Value = 61906.945
$ws.Range("K31").Value = 52507.95
$ws.Range("J31").Value = 72964.586
$ws.Range("L31").Value = 72964.586
$ws.Range("M31").Value = -52212.95
$ws.Range("N31").Value = -73554.586
$ws.Range("H34").Value = 61906.945
$ws.Range("M34").Value = -52305.95
$ws.Range("K34").Value = 52507.95
$ws.Range("I34").Value = 52507.95
$ws.Range("L34").Value = 72964.586
$ws.Range("J34").Value = 72964.586
$ws.Range("N34").Value = -73368.586
$ws.Range("I58").Value = 24391740
$ws.Range("H58").Value = 13515165
$ws.Range("M58").Value = -24391537
$ws.Range("K58").Value = 24391740
$ws.Range("H99").Value = 2150
$ws.Range("J99").Value = 2500
$ws.Range("N99").Value = -5496
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 2500
$ws.Range("I99").Value = 1800
$ws.Range("M99").Value = -302
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2930
$ws.Range("J126").Value = 2500
$ws.Range("H126").Value = 2150
$ws.Range("N126").Value = -12440
$ws.Range("I126").Value = 1800
$ws.Range("J132").Value = 127998.125
$ws.Range("L132").Value = 383994.375
$ws.Range("H132").Value = 35533.5
$ws.Range("N132").Value = -389054.375
$ws.Range("I132").Value = 1910
$ws.Range("M132").Value = -3200
$ws.Range("K132").Value = 5730
$ws.Range("M136").Value = -73172670
$ws.Range("H136").Value = 13515165
$ws.Range("K136").Value = 73175220
$ws.Range("I136").Value = 24391740

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K92").Value = 4000.5
$ws.Range("I92").Value = 1333.5
$ws.Range("M92").Value = -2752.5
$ws.Range("H92").Value = 1333.5
$ws.Range("L131").Value = 2817.41181
$ws.Range("N131").Value = -12897.41181
$ws.Range("J131").Value = 939.1372699999999
$ws.Range("H131").Value = 864.0164

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N102").Value = -5076.3636
$ws.Range("J102").Value = 1832.3636
$ws.Range("K102").Value = 1896.15
$ws.Range("M102").Value = -274.1500000000001
$ws.Range("H102").Value = 1873.5161
$ws.Range("L102").Value = 1832.3636
$ws.Range("I102").Value = 1896.15
$ws.Range("H122").Value = 1963.1538
$ws.Range("I122").Value = 1635.6451
$ws.Range("J122").Value = 3232.25
$ws.Range("M122").Value = -2456.9353
$ws.Range("N122").Value = -14596.75
$ws.Range("K122").Value = 4906.9353
$ws.Range("L122").Value = 9696.75
$ws.Range("K126").Value = 4938.6666
$ws.Range("M126").Value = -2468.6666
$ws.Range("H126").Value = 1781.6
$ws.Range("I126").Value = 1646.2222
$ws.Range("J132").Value = 257503
$ws.Range("L132").Value = 772509
$ws.Range("H132").Value = 79944.766
$ws.Range("N132").Value = -777569
$ws.Range("I132").Value = 47661.453
$ws.Range("M132").Value = -140454.359
$ws.Range("K132").Value = 142984.359
$ws.Range("J133").Value = 57920
$ws.Range("N133").Value = -68040
$ws.Range("H133").Value = 57920
$ws.Range("L133").Value = 57920

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 41669450
$ws.Range("I7").Value = 41669450
$ws.Range("M7").Value = -41669338
$ws.Range("H7").Value = 41669450
$ws.Range("K9").Value = 423.33334
$ws.Range("I9").Value = 423.33334
$ws.Range("H9").Value = 1057.5454
$ws.Range("M9").Value = -199.33334
$ws.Range("K40").Value = 2489.6
$ws.Range("M40").Value = -2353.6
$ws.Range("H40").Value = 2681.3125
$ws.Range("I40").Value = 2489.6
$ws.Range("L74").Value = 36217
$ws.Range("J74").Value = 36217
$ws.Range("N74").Value = -38213
$ws.Range("H74").Value = 2025769.6
$ws.Range("H77").Value = 2025769.6
$ws.Range("J77").Value = 36217
$ws.Range("N77").Value = -118635
$ws.Range("L77").Value = 108651
$ws.Range("H119").Value = 22950
$ws.Range("J119").Value = 22950
$ws.Range("L119").Value = 22950
$ws.Range("N119").Value = -32626
$ws.Range("H122").Value = 3578.4324
$ws.Range("I122").Value = 3459.6365
$ws.Range("J122").Value = 3752.6667
$ws.Range("M122").Value = -7928.9095
$ws.Range("N122").Value = -16158.0001
$ws.Range("K122").Value = 10378.9095
$ws.Range("L122").Value = 11258.0001
$ws.Range("K126").Value = 125008350
$ws.Range("M126").Value = -125005880
$ws.Range("H126").Value = 41669450
$ws.Range("I126").Value = 41669450

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H122").Value = 1794.8846
$ws.Range("I122").Value = 1518.1333
$ws.Range("J122").Value = 2172.2727
$ws.Range("M122").Value = -2104.3999
$ws.Range("N122").Value = -11416.8181
$ws.Range("K122").Value = 4554.3999
$ws.Range("L122").Value = 6516.8181
$ws.Range("K126").Value = 2167.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = 302.5
$ws.Range("J126").Value = 5000
$ws.Range("H126").Value = 1578
$ws.Range("N126").Value = -19940
$ws.Range("I126").Value = 722.5
$ws.Range("N136").Value = -252510
$ws.Range("M136").Value = -64575.831
$ws.Range("L136").Value = 247410
$ws.Range("H136").Value = 35395.8
$ws.Range("J136").Value = 82470
$ws.Range("K136").Value = 67125.831
$ws.Range("I136").Value = 22375.277
